$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value while forcing text storage (prevents Excel from
# auto-converting numeric-looking strings like "1.00" or "47.10" into
# real numbers, which would also touch cell number formatting).
function Set-TextValue {
    param($CellRef, $Val)
    $helper = $ws.Range('ZZ1')
    $helper.NumberFormat = '@'
    $helper.Value = $Val
    $helper.Copy()
    $ws.Range($CellRef).PasteSpecial(-4163)
}

$ws.Range('D2').Value = '65.382.02'
$ws.Range('E2').Value = '  -4.43%  '
$ws.Range('D3').Value = '3.250.89'
$ws.Range('E3').Value = '  -5.41%  '
$ws.Range('E4').Value = '  +0.01%  '
Set-TextValue 'D5' '555.41'
$ws.Range('E5').Value = '  -2.80%  '
Set-TextValue 'D6' '180.38'
$ws.Range('E6').Value = '  -4.57%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -1.97%  '
$ws.Range('D9').Value = '3.243.30'
$ws.Range('E9').Value = '  -5.25%  '
Set-TextValue 'D10' '0.184'
$ws.Range('E10').Value = '  -8.27%  '
Set-TextValue 'D11' '0.584'
$ws.Range('E11').Value = '  -4.35%  '
Set-TextValue 'D12' '47.10'
$ws.Range('E12').Value = '  -7.19%  '
Set-TextValue 'D13' '0.0000264'
$ws.Range('E13').Value = '  -6.27%  '
Set-TextValue 'D14' '633.55'
$ws.Range('E14').Value = '  +0.42%  '
Set-TextValue 'D15' '8.53'
$ws.Range('E15').Value = '  -5.25%  '
$ws.Range('D16').Value = '3.784.02'
$ws.Range('E16').Value = '  -5.10%  '
$ws.Range('D17').Value = '65.412.45'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('E18').Value = '  -3.16%  '
Set-TextValue 'D19' '17.63'
$ws.Range('E19').Value = '  -1.92%  '
$ws.Range('D20').Value = '3.259.61'
$ws.Range('E20').Value = '  -5.79%  '
Set-TextValue 'D21' '11.31'
$ws.Range('E21').Value = '  -7.19%  '
Set-TextValue 'D22' '0.898'
$ws.Range('E22').Value = '  -3.66%  '
Set-TextValue 'D23' '17.73'
$ws.Range('E23').Value = '  +0.28%  '
Set-TextValue 'D24' '104.76'
$ws.Range('E24').Value = '  +6.43%  '
Set-TextValue 'D25' '4.92'
$ws.Range('E25').Value = '  -7.42%  '
Set-TextValue 'D26' '3.97'
$ws.Range('E26').Value = '  -6.09%  '
Set-TextValue 'D27' '2.65'
$ws.Range('E27').Value = '  -5.74%  '
Set-TextValue 'D28' '9.46'
$ws.Range('E28').Value = '  -3.01%  '
Set-TextValue 'D29' '8.67'
$ws.Range('E29').Value = '  -4.65%  '
Set-TextValue 'D30' '30.18'
$ws.Range('E30').Value = '  -5.88%  '
Set-TextValue 'D31' '4.01'
$ws.Range('E31').Value = '  -2.29%  '
Set-TextValue 'D32' '6.27'
$ws.Range('E32').Value = '  -5.52%  '
Set-TextValue 'D33' '10.98'
$ws.Range('E33').Value = '  -4.33%  '
Set-TextValue 'D34' '545.35'
$ws.Range('E34').Value = '  +9.02%  '
$ws.Range('E35').Value = '  -2.82%  '
Set-TextValue 'D36' '0.999'
$ws.Range('E36').Value = '  -0.02%  '
Set-TextValue 'D37' '56.88'
$ws.Range('E37').Value = '  -5.98%  '
$ws.Range('D38').Value = '3.583.67'
$ws.Range('E38').Value = '  -1.11%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D39' '3.39'
$ws.Range('E39').Value = '  -1.54%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0715'
$ws.Range('E40').Value = '  -7.88%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D41' '2.71'
$ws.Range('E41').Value = '  -5.34%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 'D42' '0.129'
$ws.Range('E42').Value = '  -2.44%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
Set-TextValue 'D43' '3.37'
$ws.Range('E43').Value = '  -2.31%  '
Set-TextValue 'D44' '31.83'
$ws.Range('E44').Value = '  -6.18%  '
Set-TextValue 'D45' '3.31'
$ws.Range('E45').Value = '  -0.67%  '
Set-TextValue 'D46' '0.333'
$ws.Range('E46').Value = '  -8.41%  '
Set-TextValue 'D47' '0.0414'
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D48' '0.128'
$ws.Range('E48').Value = '  -3.26%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue 'D49' '2.59'
$ws.Range('E49').Value = '  -6.59%  '
$ws.Range('E50').Value = '  -0.12%  '
Set-TextValue 'D51' '1.23'
$ws.Range('E51').Value = '  +1.54%  '

# Clean up helper cell
$ws.Range('ZZ1').Clear()
$excel.CutCopyMode = $false
